$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (__str__ test) keeps its own wording, but the shared string it uses
# now absorbs the wording that used to live in G12 (row 12, calculate_area
# test), since that row is repurposed to show the area's expected result.
$ws.Range("G11").Value = "The shape color is red.This rectangle has four sides with the lengths of 5, 6, 5 and 6 centimeters."

# Row 13 (calculate_perimeter) previously had no expected-result text; give
# it one, matching the bold styling used by the other "Expected Result"
# cells in this column.
$ws.Range("G13").Value = "Perimeter = 22"
$ws.Range("G13").Font.Bold = $true

# Row 12 (calculate_area) now shows its own expected-result text.
$ws.Range("G12").Value = "Area = 30"

# Reflect the saved selection/scroll position.
$ws.Range("G12").Select()
$excel.ActiveWindow.ScrollRow = 9
